# "added 4wk low sales check"
# Updates the per-week forecast figures (MyForecast, Inventory Coverage,
# Seasonality Index) on "Forecast Comparison" and the roll-up totals on
# "Summary" to reflect the refreshed forecast run.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison ---------------------------------------------------

# Week W10 (row 2)
$ws1.Range("D2").Value = 64
$ws1.Range("H2").Value = 0.29
$ws1.Range("L2").Value = 1

# Week W11 (row 3)
$ws1.Range("D3").Value = 64
$ws1.Range("L3").Value = 0.98

# Week W12 (row 4)
$ws1.Range("D4").Value = 62
$ws1.Range("L4").Value = 1.1

# Week W13 (row 5)
$ws1.Range("D5").Value = 61
$ws1.Range("L5").Value = 0.89

# Week W14 (row 6)
$ws1.Range("D6").Value = 61
$ws1.Range("L6").Value = 1.13

# Week W15 (row 7)
$ws1.Range("D7").Value = 61
$ws1.Range("L7").Value = 0.89

# Week W16 (row 8)
$ws1.Range("D8").Value = 60
$ws1.Range("L8").Value = 0.8

# Week W17 (row 9)
$ws1.Range("L9").Value = 1.11

# Week W18 (row 10)
$ws1.Range("D10").Value = 59
$ws1.Range("L10").Value = 0.9

# Week W19 (row 11)
$ws1.Range("L11").Value = 1.05

# Week W20 (row 12)
$ws1.Range("L12").Value = 1.03

# Week W21 (row 13)
$ws1.Range("D13").Value = 57
$ws1.Range("L13").Value = 0.99

# Week W22 (row 14)
$ws1.Range("L14").Value = 0.89

# Week W23 (row 15)
$ws1.Range("D15").Value = 56
$ws1.Range("L15").Value = 0.99

# Week W24 (row 16)
$ws1.Range("D16").Value = 57
$ws1.Range("L16").Value = 0.95

# Week W25 (row 17)
$ws1.Range("D17").Value = 56
$ws1.Range("L17").Value = 0.93

# --- Summary -----------------------------------------------------------
# These cells store their numbers as text in the workbook, so a leading
# apostrophe is used to keep them as text instead of auto-converting to a
# numeric value.

$ws2.Range("B9").Value  = "'957"   # Total Forecast (16 Weeks)
$ws2.Range("B10").Value = "'497"   # Total Forecast (8 Weeks)
$ws2.Range("B11").Value = "'254"   # Total Forecast (4 Weeks)
$ws2.Range("B12").Value = "'65"    # Max Forecast
$ws2.Range("B14").Value = "'56"    # Min Forecast
